# SQL Phasing Data Set - add "Global Transfer" assignment (column T) and a new
# aggregated "BU" row, per commit: "Ajout du Global Transfer dans le jeu de
# donnee pour PNE_PHASE_01".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# 1) Insert a new row above the current row 3 ("Legal Employer"/CASA) so the
#    table gains a "BU" summary row. Inserting copies formatting from the row
#    above it (row 2, style index 4), which is exactly what the target file
#    uses for the new row 3.
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).Insert()
Write-Host "row inserted"

$ws.Range("A3").Value = "BU"
foreach ($col in @("B", "D", "F", "H", "J", "L", "N", "P", "R")) {
    $ws.Range($col + "3").Value = "CASA ES"
}
Write-Host "row3 filled"

# ---------------------------------------------------------------------------
# 2) The row that used to be row 3 ("Legal Employer" / "CREDIT AGRICOLE
#    S.A.") is now row 4. In the target workbook it no longer uses the small
#    grey Segoe-UI font style - it now matches the plain centred style used
#    by rows such as the old "Salary"/"Location" rows. Copy that format over
#    (cell by cell, so the untouched separator columns C/E/G/... stay empty).
# ---------------------------------------------------------------------------
$ws.Range("A10").Copy()
foreach ($col in @("A", "B", "D", "F", "H", "J", "L", "N", "P", "R")) {
    $ws.Range($col + "4").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0
Write-Host "row4 restyled"

# ---------------------------------------------------------------------------
# 3) Add the new "Global Transfer" assignment in column T.
#    Shared-string order matters: type T4 before T3, and T8 before T7, so
#    the new strings land in the same uniqueCount order as the target file.
# ---------------------------------------------------------------------------
$ws.Range("T1").Value = "PNE_PHASE_01"
$ws.Range("T2").Value = 44849
$ws.Range("T4").Value = "CREDIT AGRICOLE CIB France"
$ws.Range("T3").Value = "CACIB"
$ws.Range("T5").Value = "CDI"
$ws.Range("T6").Value = "CASA ES"
$ws.Range("T8").Value = "CACEIS Bank"
$ws.Range("T7").Value = "Trader"
$ws.Range("T9").Value = 1
$ws.Range("T10").Value = 70000
$ws.Range("T11").Value = "Location_001"
Write-Host "column T filled"

# Give T2 the same date format as the other period cells on that row.
$ws.Range("B2").Copy()
$ws.Range("T2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match styles on rows 6/8/9 (the "customFormat" rows) for the new column T
# and its separator column S, copying from column R/Q of the same row.
$ws.Range("R6").Copy()
$ws.Range("T6").PasteSpecial(-4122)
$ws.Range("Q6").Copy()
$ws.Range("S6").PasteSpecial(-4122)

$ws.Range("R8").Copy()
$ws.Range("T8").PasteSpecial(-4122)
$ws.Range("Q8").Copy()
$ws.Range("S8").PasteSpecial(-4122)

$ws.Range("R9").Copy()
$ws.Range("T9").PasteSpecial(-4122)
$ws.Range("Q9").Copy()
$ws.Range("S9").PasteSpecial(-4122)
$excel.CutCopyMode = 0
Write-Host "column T restyled"

# ---------------------------------------------------------------------------
# 4) Column widths: a narrow separator column S and a bestFit-style data
#    column T, matching the existing alternating layout.
# ---------------------------------------------------------------------------
$ws.Range("T1").EntireColumn.ColumnWidth = 14.1666666666667
$ws.Range("S1").EntireColumn.ColumnWidth = 0.53
Write-Host "columns sized"

# ---------------------------------------------------------------------------
# 5) Selection / active cell, matching the saved state in the target file.
# ---------------------------------------------------------------------------
$ws.Range("T11").Select()

Write-Host "done"
